$d = $word.ActiveDocument

# Titles for the new "index" paragraphs that replace the old empty
# bookmark paragraph (paragraph 2). The last two entries are the
# trailing blank paragraphs added at the end.
$titles = @(
    "Que es .Net",
    "Lenguajes que usan .Net",
    "Tipos de Datos",
    "Clases",
    "Templates (Plantillas)",
    "Listas",
    "Interfaces Gráficas I",
    "Interfaces Gráficas II",
    "Conexión a Bases de Datos (SqlServer)",
    "Conexión a Bases de Datos (MySql)",
    "Conexión a Bases de Datos (PostgreSql)",
    "",
    ""
)

# Paragraph 2 is the (currently empty) paragraph carrying the
# "_GoBack" bookmark and the centered/bold/size-50 formatting.
# Strip the centering so the new index lines are left aligned, then
# turn it into the first index entry ("Que es .Net").
$p2 = $d.Paragraphs(2)
$p2.Alignment = 0
$p2.Range.InsertBefore($titles[0])
$p2b = $d.Paragraphs(2)
$p2b.Range.Font.Bold = 1
$p2b.Range.Font.Size = 25
$p2b.Range.Font.SizeBi = 25

# Remove the original hidden "_GoBack" bookmark; it will be re-added
# on the last populated index entry further down.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Add the remaining paragraphs, one per remaining title (including the
# two trailing empty ones), inheriting bold/size-50 formatting without
# centering.
for ($i = 1; $i -lt $titles.Length; $i++) {
    $prevIndex = 1 + $i
    $prev = $d.Paragraphs($prevIndex)
    $prev.Range.InsertParagraphAfter()

    $curIndex = $prevIndex + 1
    $cur = $d.Paragraphs($curIndex)
    $cur.Alignment = 0

    $text = $titles[$i]
    if ($text -ne "") {
        $cur.Range.InsertBefore($text)
        $cur2 = $d.Paragraphs($curIndex)
        $cur2.Range.Font.Bold = 1
        $cur2.Range.Font.Size = 25
        $cur2.Range.Font.SizeBi = 25
    }
}

# Re-create the "_GoBack" bookmark on the last non-empty index entry
# ("Conexión a Bases de Datos (PostgreSql)").
$pgIndex = 2 + 10
$pgPara = $d.Paragraphs($pgIndex)
$endRange = $pgPara.Range
$endRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRange)

Write-Output ("PARAGRAPHS=" + $d.Paragraphs.Count)
Write-Output ("TEXT=[" + $d.Content.Text + "]")
